$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing daily stats (AgTests / AgPosit corrections) ---
$ws.Range("F322").Value = 107320
$ws.Range("G322").Value = 2306
$ws.Range("F323").Value = 212233
$ws.Range("G323").Value = 3158
$ws.Range("F324").Value = 233051
$ws.Range("G324").Value = 2663
$ws.Range("F325").Value = 750784
$ws.Range("G325").Value = 6323
$ws.Range("F326").Value = 426844
$ws.Range("G326").Value = 3755
$ws.Range("F327").Value = 238962
$ws.Range("G327").Value = 2885
$ws.Range("F328").Value = 180445
$ws.Range("G328").Value = 2645
$ws.Range("F329").Value = 88936
$ws.Range("G329").Value = 1802
$ws.Range("F331").Value = 150584
$ws.Range("G331").Value = 2587
$ws.Range("F332").Value = 422944
$ws.Range("G332").Value = 4109
$ws.Range("F333").Value = 258608
$ws.Range("F334").Value = 202288
$ws.Range("G334").Value = 3380
$ws.Range("F335").Value = 129634
$ws.Range("G335").Value = 2903
$ws.Range("F336").Value = 100226
$ws.Range("F337").Value = 101962
$ws.Range("F338").Value = 218482
$ws.Range("G338").Value = 3075
$ws.Range("F339").Value = 625041
$ws.Range("G339").Value = 5361
$ws.Range("F340").Value = 372465
$ws.Range("G340").Value = 3192
$ws.Range("F341").Value = 295228
$ws.Range("G341").Value = 3642
$ws.Range("F342").Value = 184669
$ws.Range("G342").Value = 3118
$ws.Range("F343").Value = 122386
$ws.Range("G343").Value = 2754

# --- Append new row for 2021-02-12 (pi 12. 02. 2021) ---
$ws.Range("A344").Value = 44238
$ws.Range("A344").NumberFormat = "yyyy-mm-dd"
$ws.Range("B344").Value = 273904
$ws.Range("C344").Value = 10402
$ws.Range("D344").Value = 2431
$ws.Range("E344").Value = 5733
$ws.Range("F344").Value = 101684
$ws.Range("G344").Value = 1917
